$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new top data row for 2022-Q4 and push the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make room: copy row 6 (last existing data row, 2021-Q3) down into the new
# row 7 so the row-number cell (column A) keeps its border/bold style.
$summary.Range("A6:D6").Copy()
$summary.Range("A7:D7").PasteSpecial(-4122)

# Shift the quarter rows down by one (values only - formatting already in
# place on rows 2-6).
$summary.Range("B7").Value = "2021-Q3"
$summary.Range("C7").Value = 7
$summary.Range("D7").Value = 0.85

$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 15
$summary.Range("D6").Value = 6.34

$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 7
$summary.Range("D5").Value = 2.28

$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 8
$summary.Range("D4").Value = 1.74

$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 6
$summary.Range("D3").Value = 1.03

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 0.79

# Row-number helper column (A) just keeps counting up 0..5.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------------
# 2) Add the brand-new "2022-Q4" sheet right before "2022-Q3". Duplicating
#    the existing "2022-Q3" sheet (instead of Worksheets.Add()) means the new
#    sheet inherits its styles/number formats exactly, with no extra style
#    bookkeeping to redo by hand.
# ---------------------------------------------------------------------------
$existingQ3 = $wb.Worksheets.Item("2022-Q3")
$existingQ3.Copy($existingQ3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template only has 6 data rows (rows 2-7); the new data needs 8, so
# clone the last data row's formatting down into rows 8 and 9.
$q4.Range("A7:H7").Copy()
$q4.Range("A8:H8").PasteSpecial(-4122)
$q4.Range("A9:H9").PasteSpecial(-4122)

$rows = @(
    @(0, "001766", "上投摩根医疗健康股票A", "9.04", "84.27", "3.89", "0.3517", 8),
    @(1, "630010", "华商价值精选混合", "4.37", "87.25", "3.10", "0.1355", 9),
    @(2, "002666", "前海开源沪港深创新成长灵活配置混合A", "4.54", "67.11", "2.51", "0.1140", 8),
    @(3, "002667", "前海开源沪港深创新成长灵活配置混合C", "2.14", "67.11", "2.51", "0.0537", 8),
    @(4, "014932", "上投摩根医疗健康股票C", "1.32", "84.27", "3.89", "0.0513", 8),
    @(5, "588160", "南方上证科创板新材料ETF", "0.90", "98.46", "3.79", "0.0341", 7),
    @(6, "630006", "华商产业升级混合", "0.86", "88.65", "3.16", "0.0272", 9),
    @(7, "588010", "博时上证科创板新材料ETF", "0.53", "98.90", "3.82", "0.0202", 7)
)

# Columns B-G are stored as TEXT in the source data (fund code keeps its
# leading zeros, the numeric-looking figures are plain strings too) while
# column A (row counter) and H (position rank) are real numbers. Force the
# "@" text format on B:G before assigning so numeric-looking strings like
# "9.04" or "001766" are not auto-coerced into numbers.
$q4.Range("B2:G9").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $rowNum = $r + 2
    $q4.Cells.Item($rowNum, 1).Value = $row[0]
    $q4.Cells.Item($rowNum, 2).Value = $row[1]
    $q4.Cells.Item($rowNum, 3).Value = $row[2]
    $q4.Cells.Item($rowNum, 4).Value = $row[3]
    $q4.Cells.Item($rowNum, 5).Value = $row[4]
    $q4.Cells.Item($rowNum, 6).Value = $row[5]
    $q4.Cells.Item($rowNum, 7).Value = $row[6]
    $q4.Cells.Item($rowNum, 8).Value = $row[7]
}
